# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.950.11"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.843.82"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'232.44"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'41.27"
$ws.Range("E8").Value = "  +6.83%  "
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").Value = "'0.0692"
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "'0.0985"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "2.112.24"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'11.38"
$ws.Range("E13").Value = "  +5.36%  "
$ws.Range("D14").Value = "1.841.42"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'0.671"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "'4.66"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "34.951.93"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'69.93"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("D20").Value = "'240.39"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'12.17"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "'171.89"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'17.45"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").Value = "'0.123"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("E29").Value = "  +9.42%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'0.0552"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "'3.95"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'3.90"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'1.63"
$ws.Range("E34").Value = "  +21.59%  "
$ws.Range("D35").Value = "'1.95"
$ws.Range("E35").Value = "  +11.71%  "
$ws.Range("D36").Value = "'0.744"
$ws.Range("E36").Value = "  +9.23%  "
$ws.Range("D37").Value = "'1.22"
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("E38").Value = "  +11.96%  "
$ws.Range("D39").Value = "'89.60"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "1.348.18"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "'14.51"
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").Value = "'2.27"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").Value = "'0.0530"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "'6.32"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").Value = "2.030.38"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.01"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'3.40"
$ws.Range("E50").Value = "  +15.40%  "
$ws.Range("D51").Value = "'0.0668"
$ws.Range("E51").Value = "  -2.63%  "
